$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of rows 67 and 68 (columns F:V); columns A-E are unchanged ---
$tmp67 = $ws.Range("F67:V67").Value2
$tmp68 = $ws.Range("F68:V68").Value2
$ws.Range("F68:V68").Value = $tmp67
$ws.Range("F67:V67").Value = $tmp68

# --- Swap the contents of rows 74 and 75 (columns F:V); columns A-E are unchanged ---
$tmp74 = $ws.Range("F74:V74").Value2
$tmp75 = $ws.Range("F75:V75").Value2
$ws.Range("F75:V75").Value = $tmp74
$ws.Range("F74:V74").Value = $tmp75

# --- Append new row 77 (new match: Al Ain vs Khorfakkan) ---
# Copy row 76 first so the new row inherits the same formatting/styles.
$ws.Range("A76:V76").Copy($ws.Range("A77:V77"))

$ws.Range("A77").Value = 76
$ws.Range("B77").Value = "united-arab-emirates"
$ws.Range("C77").Value = "uae-league"
$ws.Range("D77").Value = "2023-2024"
$ws.Range("E77").Value = 45280.57291666666
$ws.Range("F77").Value = "Al Ain"
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = "Khorfakkan"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 1.18
$ws.Range("K77").Value = "16/12/2023 18:13"
$ws.Range("L77").Value = 1.16
$ws.Range("M77").Value = "20/12/2023 13:38"
$ws.Range("N77").Value = 7.38
$ws.Range("O77").Value = "16/12/2023 18:13"
$ws.Range("P77").Value = 8.44
$ws.Range("Q77").Value = "20/12/2023 13:43"
$ws.Range("R77").Value = 9.81
$ws.Range("S77").Value = "16/12/2023 18:13"
$ws.Range("T77").Value = 13.48
$ws.Range("U77").Value = "20/12/2023 13:43"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-ain-khorfakkan/KrvszPf3/"
